$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "1.0000", "0.7124") are preserved verbatim instead of being
# coerced into numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.448.41"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.876.59"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "0.7124"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "241.89"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.07842"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").Value = "0.3109"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("D10").Value = "25.17"
$ws.Range("E10").Value = "  +7.36%  "
$ws.Range("D11").Value = "0.08240"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.883.06"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7273"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").Value = "5.261"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "90.88"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "29.448.68"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "5.914"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").Value = "247.43"
$ws.Range("E18").Value = "  +4.37%  "
$ws.Range("D19").Value = "0.000007868"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "13.26"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "2.121.73"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "7.970"
$ws.Range("E23").Value = "  +7.27%  "
$ws.Range("D24").Value = "0.9993"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "0.1572"
$ws.Range("E25").Value = "  +10.52%  "
$ws.Range("D26").Value = "163.68"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "9.014"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "18.28"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("D29").Value = "1.365"
$ws.Range("E29").Value = "  -3.34%  "
$ws.Range("D30").Value = "1.494"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "4.370"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").Value = "4.125"
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").Value = "0.05307"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").Value = "1.931"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").Value = "1.199"
$ws.Range("E35").Value = "  +3.30%  "
$ws.Range("D36").Value = "0.7229"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").Value = "1.246.15"
$ws.Range("E39").Value = "  +8.10%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("D42").Value = "74.08"
$ws.Range("E42").Value = "  +5.42%  "
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "103.29"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "0.5323"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "1.773"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("D48").Value = "2.926"
$ws.Range("E48").Value = "  +12.87%  "
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "0.4319"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").Value = "9.251"
$ws.Range("E51").Value = "  +1.16%  "
